# Force Sensor Diagram - update pump PSI input and refresh selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# PSI input (D7 label "PSI") drops from 2000 to 1500 -> Force (E8 = E7*E6) recalculates.
$ws.Range("E7").Value = 1500

# Move the active selection to the cell that was edited.
$ws.Range("E7").Select()
